$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '35.125.43'
$ws.Range('E2').Value = '  -0.50%  '
$ws.Range('D3').Value = '1.900.60'
$ws.Range('E3').Value = '  +0.13%  '
$ws.Range('E4').Value = '  -0.45%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '252.53'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.55%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.701'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.28%  '
$ws.Range('E7').Value = '  -0.38%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '41.65'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +3.16%  '
$ws.Range('E9').Value = '  +2.05%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '52.31'
$ws.Range('D10').Style = 'Normal'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0757'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +4.86%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0979'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.99%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '13.05'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +4.12%  '
$ws.Range('D14').Value = '2.177.11'
$ws.Range('E14').Value = '  +0.22%  '
$ws.Range('E15').Value = '  +3.34%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '4.99'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +3.25%  '
$ws.Range('D17').Value = '1.897.16'
$ws.Range('E17').Value = '  +0.14%  '
$ws.Range('D18').Value = '35.139.94'
$ws.Range('E18').Value = '  -0.32%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '73.98'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +2.04%  '
$ws.Range('E20').Value = '  +2.54%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '249.72'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +3.67%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '13.03'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +2.02%  '
$ws.Range('E23').Value = '  +2.23%  '
$ws.Range('E24').Value = '  -0.42%  '
$ws.Range('E25').Value = '  +4.78%  '
$ws.Range('E26').Value = '  -1.64%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '168.56'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.64%  '
$ws.Range('E28').Value = '  +0.21%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '18.51'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -2.56%  '
$ws.Range('E30').Value = '  -0.56%  '
$ws.Range('D31').Value = '4.128.37'
$ws.Range('E31').Value = '  -0.34%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.08'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +10.98%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.33'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +3.41%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0596'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +4.63%  '
$ws.Range('E35').Value = '  +10.69%  '
$ws.Range('E36').Value = '  +3.85%  '
$ws.Range('E37').Value = '  -0.47%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.850'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -6.75%  '
$ws.Range('E39').Value = '  +0.45%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '17.57'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +7.33%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '98.81'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +3.25%  '
$ws.Range('B42').Value = 'VeChain'
$ws.Range('C42').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0215'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +3.86%  '
$ws.Range('B43').Value = 'Kaspa'
$ws.Range('C43').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0671'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +2.25%  '
$ws.Range('E44').Value = '  +0.62%  '
$ws.Range('E45').Value = '  +1.27%  '
$ws.Range('D46').Value = '1.309.57'
$ws.Range('E46').Value = '  -3.40%  '
$ws.Range('E47').Value = '  -0.08%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.75'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.38%  '
$ws.Range('E49').Value = '  +2.05%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '12.14'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.15%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0761'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +7.97%  '
